# Add the new 근로장학생 선발 구현 row to the 이채현 worksheet and make that
# sheet/cell the active selection, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("이채현")

$ws.Range("A8").Value = "근로장학생 선발 구현"
$ws.Range("B8").Value = "성적만을 기준으로 하는 근로장학생 선발 기능을 구현함."
$ws.Range("C8").Value = "2019-05-20"
$ws.Range("D8").Value = "2019-05-20"
$ws.Range("E8").Value = "완료"
$ws.Range("F8").Value = "테스트 미실시"

$null = $ws.Activate()
$null = $ws.Range("F8").Select()
